$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append starting at row 8 (A:I)
# Columns: A=Num Nodes, B=Num Seeds, C=Timestamp (text), D=Greedy DP, E=Greedy,
#          F=Recursive DP Root, G=Rec DP No Root, H=LP, I=Bipartite LP
$rows = @(
    @(20, 10, "2021-04-01 15:14:42.328650", 11, 11, 11, 4, 11, 12),
    @(100, 10, "2021-04-01 15:27:37.709801", 34, 34, 32, 34, 34, 36),
    @(100, 10, "2021-04-01 15:28:44.328220", 42, 42, 42, 32, 42, 45),
    @(100, 10, "2021-04-01 15:28:55.951325", 37, 36, 37, 36, 37, 42),
    @(100, 10, "2021-04-01 15:29:05.953316", 27, 26, 27, 22, 27, 31),
    @(50, 5, "2021-04-08 10:28:10.365738", 11, 11, 11, 11, 11, 2),
    @(50, 5, "2021-04-08 10:30:09.387693", 30, 30, 30, 3, 30, 30),
    @(50, 5, "2021-04-08 10:30:20.690811", 17, 18, 17, 18, 18, 18),
    @(100, 10, "2021-04-08 10:31:25.854300", 45, 43, 46, 18, 46, 46),
    @(100, 10, "2021-04-08 10:32:02.904430", 46, 47, 47, 33, 47, 47),
    @(100, 10, "2021-04-08 10:33:05.418990", 47, 48, 48, 47, 48, 48)
)

$startRow = 8
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
}
